$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure percentage-like text values stay as text (not auto-converted to numbers)
$pctCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10")
foreach ($addr in $pctCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("B2").Value = "34°"
$ws.Range("C2").Value = "19°"
$ws.Range("D2").Value = "34%"
$ws.Range("E2").Value = "69%"

# Row 3
$ws.Range("C3").Value = "20°"
$ws.Range("D3").Value = "38%"
$ws.Range("E3").Value = "52%"

# Row 4
$ws.Range("E4").Value = "48%"

# Row 5
$ws.Range("B5").Value = "28°"
$ws.Range("D5").Value = "58%"
$ws.Range("E5").Value = "83%"

# Row 6
$ws.Range("B6").Value = "21°"
$ws.Range("C6").Value = "16°"
$ws.Range("D6").Value = "71%"
$ws.Range("E6").Value = "86%"
$ws.Range("F6").Value = "8 de 11"

# Row 7
$ws.Range("B7").Value = "26°"
$ws.Range("D7").Value = "67%"
$ws.Range("E7").Value = "85%"

# Row 8
$ws.Range("B8").Value = "32°"
$ws.Range("C8").Value = "19°"
$ws.Range("D8").Value = "52%"
$ws.Range("E8").Value = "66%"

# Row 9
$ws.Range("B9").Value = "35°"
$ws.Range("C9").Value = "21°"
$ws.Range("D9").Value = "34%"
$ws.Range("E9").Value = "51%"

# Row 10
$ws.Range("B10").Value = "32°"
$ws.Range("C10").Value = "19°"
$ws.Range("D10").Value = "40%"
$ws.Range("E10").Value = "68%"
$ws.Range("F10").Value = "10 de 11"

# Remove row 11 entirely (shift remaining rows up, update dimension)
$ws.Rows("11:11").Delete()
